$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C12").Value = 0.158
$ws.Range("C13").Value = -2.045
$ws.Range("C39").Value = -0.53
$ws.Range("C45").Value = -1.407
$ws.Range("C57").Value = -1.84
$ws.Range("C66").Value = -2.23
$ws.Range("C76").Value = -2.807
$ws.Range("C82").Value = -1.738
$ws.Range("C90").Value = -0.797
$ws.Range("C94").Value = 2.511
$ws.Range("C97").Value = 1.12
